$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$v_A2 = @"
A
"@
$ws.Range("A2").Value = $v_A2

$v_B2 = @"
A, or a, is the first letter and the first vowel of the modern English alphabet and the ISO basic Latin alphabet.[1][2] Its name in English is a (pronounced /ˈeɪ/), plural aes.[nb 1] It is similar in shape to the Ancient Greek letter alpha, from which it derives.[3] The uppercase version consists of the two slanting sides of a triangle, crossed in the middle by a horizontal bar. The lowercase version can be written in two forms: the double-storey a and single-storey ɑ. The latter is commonly used in handwriting and fonts based on it, especially fonts intended to be read by children, and is also found in italic type.

"@
$ws.Range("B2").Value = $v_B2

$v_A3 = @"
History
"@
$ws.Range("A3").Value = $v_A3

$v_B3 = @"
The earliest certain ancestor of "A" is aleph (also written 'aleph), the first letter of the Phoenician alphabet,[4] which consisted entirely of consonants (for that reason, it is also called an abjad to distinguish it from a true alphabet). In turn, the ancestor of aleph may have been a pictogram of an ox head in proto-Sinaitic script[5] influenced by Egyptian hieroglyphs, styled as a triangular head with two horns extended.
When the ancient Greeks adopted the alphabet, they had no use for a letter to represent the glottal stop—the consonant sound that the letter denoted in Phoenician and other Semitic languages, and that was the first phoneme of the Phoenician pronunciation of the letter—so they used their version of the sign to represent the vowel /a/, and called it by the similar name of alpha. In the earliest Greek inscriptions after the Greek Dark Ages, dating to the 8th century BC, the letter rests upon its side, but in the Greek alphabet of later times it generally resembles the modern capital letter, although many local varieties can be distinguished by the shortening of one leg, or by the angle at which the cross line is set.
The Etruscans brought the Greek alphabet to their civilization in the Italian Peninsula and left the letter unchanged. The Romans later adopted the Etruscan alphabet to write the Latin language, and the resulting letter was preserved in the Latin alphabet that would come to be used to write many languages, including English.
Typographic variantsDuring Roman times, there were many variant forms of the letter "A". First was the monumental or lapidary style, which was used when inscribing on stone or other "permanent" media. There was also a cursive style used for everyday or utilitarian writing, which was done on more perishable surfaces. Due to the "perishable" nature of these surfaces, there are not as many examples of this style as there are of the monumental, but there are still many surviving examples of different types of cursive, such as majuscule cursive, minuscule cursive, and semicursive minuscule. Variants also existed that were intermediate between the monumental and cursive styles. The known variants include the early semi-uncial, the uncial, and the later semi-uncial.[6]
At the end of the Roman Empire (5th century AD), several variants of the cursive minuscule developed through Western Europe. Among these were the semicursive minuscule of Italy, the Merovingian script in France, the Visigothic script in Spain, and the Insular or Anglo-Irish semi-uncial or Anglo-Saxon majuscule of Great Britain. By the 9th century, the Caroline script, which was very similar to the present-day form, was the principal form used in book-making, before the advent of the printing press. This form was derived through a combining of prior forms.[6]
15th-century Italy saw the formation of the two main variants that are known today. These variants, the Italic and Roman forms, were derived from the Caroline Script version. The Italic form, also called script a, is used in most current handwriting; it consists of a circle and vertical stroke on the right ("ɑ"). This slowly developed from the fifth-century form resembling the Greek letter tau in the hands of medieval Irish and English writers.[4] The Roman form is used in most printed material; it consists of a small loop with an arc over it ("a").[6] Both derive from the majuscule (capital) form. In Greek handwriting, it was common to join the left leg and horizontal stroke into a single loop, as demonstrated by the uncial version shown. Many fonts then made the right leg vertical. In some of these, the serif that began the right leg stroke developed into an arc, resulting in the printed form, while in others it was dropped, resulting in the modern handwritten form. Graphic designers refer to the Italic and Roman forms as "single decker a" and "double decker a" respectively.
Italic type is commonly used to mark emphasis or more generally to distinguish one part of a text from the rest (set in Roman type). There are some other cases aside from italic type where script a ("ɑ"), also called Latin alpha, is used in contrast with Latin "a" (such as in the International Phonetic Alphabet).

"@
$ws.Range("B3").Value = $v_B3

$v_A4 = @"
Use in writing systems
"@
$ws.Range("A4").Value = $v_A4

$v_B4 = @"
EnglishIn modern English orthography, the letter ⟨a⟩ represents at least seven different vowel sounds:
the near-open front unrounded vowel /æ/ as in pad;
the open back unrounded vowel /ɑː/ as in father, which is closer to its original Latin and Greek sound;[5]
the diphthong /eɪ/ as in ace and major (usually when ⟨a⟩ is followed by one, or occasionally two, consonants and then another vowel letter) – this results from Middle English lengthening followed by the Great Vowel Shift;
the modified form of the above sound that occurs before ⟨r⟩, as in square and Mary;
the rounded vowel of water;
the shorter rounded vowel (not present in General American) in was and what;[4]
a schwa, in many unstressed syllables, as in about, comma, solar.The double ⟨aa⟩ sequence does not occur in native English words, but is found in some words derived from foreign languages such as Aaron and aardvark.[7] However, ⟨a⟩ occurs in many common digraphs, all with their own sound or sounds, particularly ⟨ai⟩, ⟨au⟩, ⟨aw⟩, ⟨ay⟩, ⟨ea⟩ and ⟨oa⟩.
⟨a⟩ is the third-most-commonly used letter in English (after ⟨e⟩ and ⟨t⟩),[8] and the second most common in Spanish and French. In one study, on average, about 3.68% of letters used in English texts tend to be ⟨a⟩, while the number is 6.22% in Spanish and 3.95% in French.[8]
Other languagesIn most languages that use the Latin alphabet, ⟨a⟩ denotes an open unrounded vowel, such as /a/, /ä/, or /ɑ/. An exception is Saanich, in which ⟨a⟩ (and the glyph Á) stands for a close-mid front unrounded vowel /e/.
Other systemsIn phonetic and phonemic notation:
in the International Phonetic Alphabet, ⟨a⟩ is used for the open front unrounded vowel, ⟨ä⟩ is used for the open central unrounded vowel, and ⟨ɑ⟩ is used for the open back unrounded vowel.
in X-SAMPA, ⟨a⟩ is used for the open front unrounded vowel and ⟨A⟩ is used for the open back unrounded vowel.
"@
$ws.Range("B4").Value = $v_B4

$v_A5 = @"
Other uses
"@
$ws.Range("A5").Value = $v_A5

$v_B5 = @"
In algebra, the letter a along with various other letters of the alphabet is often used to denote a variable, with various conventional meanings in different areas of mathematics. Moreover, in 1637, René Descartes "invented the convention of representing unknowns in equations by x, y, and z, and knowns by a, b, and c",[9] and this convention is still often followed, especially in elementary algebra.
In geometry, capital A, B, C etc. are used to denote segments, lines, rays, etc.[6] A capital A is also typically used as one of the letters to represent an angle in a triangle, the lowercase a representing the side opposite angle A.[5]
"A" is often used to denote something or someone of a better or more prestigious quality or status: A-, A or A+, the best grade that can be assigned by teachers for students' schoolwork; "A grade" for clean restaurants; A-list celebrities, etc. Such associations can have a motivating effect, as exposure to the letter A has been found to improve performance, when compared with other letters.[10]
"A" is used as a prefix on some words, such as asymmetry, to mean "not" or "without" (from Greek).
In English grammar, "a", and its variant "an", is an indefinite article, used to introduce noun phrases.
Finally, the letter A is used to denote size, as in a narrow size shoe,[5] or a small cup size in a brassiere.[11]

"@
$ws.Range("B5").Value = $v_B5
